$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: only the ShipmentTracking value (P2) was removed; ActualRate (Q2)
# and Result (R2) stay as they were.
$ws.Range("P2").ClearContents()

# Rows 3-25: ShipmentTracking/ActualRate/Result (P:R) all cleared out,
# leaving plain blank cells with no explicit style.
$ws.Range("P3:R25").ClearContents()

# Row 26 (the last data row): same P:R clear, but this row ends up with
# explicit border formatting on P26/Q26 and border+center on R26 - apply
# the formatting first so ClearContents leaves the cells blank-but-styled.
# (Color must be set before LineStyle, otherwise an extra transient/unused
# border style is left behind in the style table.)
$ws.Range("P26:Q26").Borders.Color = 0
$ws.Range("P26:Q26").Borders.LineStyle = 1

$ws.Range("R26").Borders.Color = 0
$ws.Range("R26").Borders.LineStyle = 1
$ws.Range("R26").HorizontalAlignment = -4108

$ws.Range("P26:R26").ClearContents()
